$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / masthead text updates ---
$ws.Range("M6").Value = "Jessica S. Tisch"
$ws.Range("A8").Value = "Volume 31   Number  48"
$ws.Range("C9").Value = "Report Covering the Week  11/25/2024  Through  12/1/2024"

# --- Column width tweaks (I and J narrow slightly, matches bestFit recalculation) ---
$ws.Columns.Item(9).ColumnWidth = 5.43
$ws.Columns.Item(10).ColumnWidth = 5.43

# --- Donor cells used to clone style+placeholder text exactly for cells that switch
#     between a numeric value and the '0'/'***.* ' placeholder text via Copy (keeps the
#     original style index intact, matching how these placeholder cells are authored). ---
$zeroDonor = $ws.Cells.Item(15, 3)    # "0" placeholder text, normal-row style
$naDonor   = $ws.Cells.Item(15, 5)    # "***.*" placeholder text, normal-row style
$intDonor  = $ws.Cells.Item(15, 9)    # plain integer cell, normal-row style
$pctDonor  = $ws.Cells.Item(15, 11)   # plain percent cell, normal-row style

# --- Crime-complaint table updates (rows 15-28) ---
# Row 15
$ws.Cells.Item(15, 6).Value = 1

# Row 16
$ws.Cells.Item(16, 4).Value = 1
$ws.Cells.Item(16, 6).Value = 2
$ws.Cells.Item(16, 7).Value = 5
$ws.Cells.Item(16, 8).Value = -60
$ws.Cells.Item(16, 10).Value = 116
$ws.Cells.Item(16, 11).Value = -25.862068965517
$ws.Cells.Item(16, 12).Value = -35.338345864661
$ws.Cells.Item(16, 13).Value = -28.333333333333
$ws.Cells.Item(16, 14).Value = -84.778761061946

# Row 17
$ws.Cells.Item(17, 4).Value = 4
$ws.Cells.Item(17, 5).Value = -75
$ws.Cells.Item(17, 6).Value = 3
$ws.Cells.Item(17, 8).Value = -66.666666666666
$ws.Cells.Item(17, 9).Value = 121
$ws.Cells.Item(17, 10).Value = 108
$ws.Cells.Item(17, 11).Value = 12.037037037037
$ws.Cells.Item(17, 12).Value = -20.915032679738
$ws.Cells.Item(17, 13).Value = 75.362318840579
$ws.Cells.Item(17, 14).Value = -50.612244897959

# Row 18
$ws.Cells.Item(18, 3).Value = 2
$ws.Cells.Item(18, 4).Value = 1
$ws.Cells.Item(18, 5).Value = 100
$ws.Cells.Item(18, 6).Value = 8
$ws.Cells.Item(18, 7).Value = 13
$ws.Cells.Item(18, 8).Value = -38.461538461538
$ws.Cells.Item(18, 9).Value = 137
$ws.Cells.Item(18, 10).Value = 162
$ws.Cells.Item(18, 11).Value = -15.432098765432
$ws.Cells.Item(18, 12).Value = -43.852459016393
$ws.Cells.Item(18, 13).Value = -38.009049773755
$ws.Cells.Item(18, 14).Value = -86.977186311787

# Row 19
$ws.Cells.Item(19, 3).Value = 14
$ws.Cells.Item(19, 4).Value = 12
$ws.Cells.Item(19, 5).Value = 16.666666666666
$ws.Cells.Item(19, 6).Value = 54
$ws.Cells.Item(19, 7).Value = 51
$ws.Cells.Item(19, 8).Value = 5.882352941176
$ws.Cells.Item(19, 9).Value = 622
$ws.Cells.Item(19, 10).Value = 630
$ws.Cells.Item(19, 11).Value = -1.269841269841
$ws.Cells.Item(19, 12).Value = 10.283687943262
$ws.Cells.Item(19, 13).Value = 122.939068100358
$ws.Cells.Item(19, 14).Value = 82.941176470588

# Row 20
$ws.Cells.Item(20, 3).Value = 2
$ws.Cells.Item(20, 4).Value = 3
$ws.Cells.Item(20, 5).Value = -33.333333333333
$ws.Cells.Item(20, 6).Value = 5
$ws.Cells.Item(20, 8).Value = -61.538461538461
$ws.Cells.Item(20, 9).Value = 106
$ws.Cells.Item(20, 10).Value = 147
$ws.Cells.Item(20, 11).Value = -27.891156462585
$ws.Cells.Item(20, 12).Value = -38.011695906432
$ws.Cells.Item(20, 13).Value = -22.058823529411
$ws.Cells.Item(20, 14).Value = -87.290167865707

# Row 21
$ws.Cells.Item(21, 6).Value = 73
$ws.Cells.Item(21, 7).Value = 91
$ws.Cells.Item(21, 8).Value = -19.780219780219
$ws.Cells.Item(21, 9).Value = 1082
$ws.Cells.Item(21, 10).Value = 1169
$ws.Cells.Item(21, 11).Value = -7.442258340461
$ws.Cells.Item(21, 12).Value = -15.402658326817
$ws.Cells.Item(21, 13).Value = 30.992736077481
$ws.Cells.Item(21, 14).Value = -64.547837483617

# Row 22
$zeroDonor.Copy($ws.Cells.Item(22, 3))
$intDonor.Copy($ws.Cells.Item(22, 4))
$ws.Cells.Item(22, 4).Value = 2
$pctDonor.Copy($ws.Cells.Item(22, 5))
$ws.Cells.Item(22, 5).Value = -100
$ws.Cells.Item(22, 7).Value = 4
$ws.Cells.Item(22, 8).Value = -75
$ws.Cells.Item(22, 10).Value = 14
$ws.Cells.Item(22, 11).Value = -7.142857142857
$ws.Cells.Item(22, 12).Value = 0

# Row 23
$zeroDonor.Copy($ws.Cells.Item(23, 6))
$ws.Cells.Item(23, 12).Value = -15.625

# Row 24
$ws.Cells.Item(24, 3).Value = 13
$ws.Cells.Item(24, 4).Value = 20
$ws.Cells.Item(24, 5).Value = -35
$ws.Cells.Item(24, 6).Value = 80
$ws.Cells.Item(24, 7).Value = 68
$ws.Cells.Item(24, 8).Value = 17.647058823529
$ws.Cells.Item(24, 9).Value = 916
$ws.Cells.Item(24, 10).Value = 833
$ws.Cells.Item(24, 11).Value = 9.963985594237
$ws.Cells.Item(24, 12).Value = -3.680336487907
$ws.Cells.Item(24, 13).Value = 61.837455830388

# Row 25
$ws.Cells.Item(25, 3).Value = 13
$ws.Cells.Item(25, 4).Value = 7
$ws.Cells.Item(25, 5).Value = 85.714285714285
$ws.Cells.Item(25, 6).Value = 55
$ws.Cells.Item(25, 7).Value = 42
$ws.Cells.Item(25, 8).Value = 30.952380952381
$ws.Cells.Item(25, 9).Value = 579
$ws.Cells.Item(25, 10).Value = 457
$ws.Cells.Item(25, 11).Value = 26.695842450765
$ws.Cells.Item(25, 12).Value = 4.136690647482

# Row 26
$ws.Cells.Item(26, 3).Value = 7
$ws.Cells.Item(26, 5).Value = 75
$ws.Cells.Item(26, 6).Value = 22
$ws.Cells.Item(26, 7).Value = 19
$ws.Cells.Item(26, 8).Value = 15.78947368421
$ws.Cells.Item(26, 9).Value = 260
$ws.Cells.Item(26, 10).Value = 225
$ws.Cells.Item(26, 11).Value = 15.555555555555
$ws.Cells.Item(26, 12).Value = 1.5625
$ws.Cells.Item(26, 13).Value = 34.715025906735

# Row 27
$ws.Cells.Item(27, 6).Value = 1

# Row 28
$intDonor.Copy($ws.Cells.Item(28, 3))
$ws.Cells.Item(28, 3).Value = 1
$ws.Cells.Item(28, 6).Value = 5
$ws.Cells.Item(28, 7).Value = 1
$ws.Cells.Item(28, 8).Value = 400
$ws.Cells.Item(28, 9).Value = 40
$ws.Cells.Item(28, 11).Value = -4.761904761904
$ws.Cells.Item(28, 12).Value = 42.857142857142

